# Bug for addition wealth class fixed
# - test for reduction of wealth class
#
# The "wealth class" figures in column B of the Main sheet were adjusted
# (some increased, to validate the addition-bug fix; one decreased, to add
# coverage for the reduction case), and the "Records Banks Column" letter
# (stored as text in B9) moved from "G" to "K". The last touched/active
# cell in the sheet is left on B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wealth Class in Allocation Row: 12 -> 22 (addition case)
$ws.Range("B2").Value = 22

# Wealth Class in Cash Flow Row: 16 -> 20 (addition case)
$ws.Range("B3").Value = 20

# Wealth Row: 7 -> 9 (reduction-of-wealth-class regression test)
$ws.Range("B5").Value = 9

# Records Banks Column letter: "G" -> "K"
$ws.Range("B9").Value = "K"

# Leave the selection on B3, matching the cell last edited/reviewed
[void]$ws.Range("B3").Select()

# Cell style "Normal" renamed to "Standard" in the source workbook's
# styles.xml; attempt the rename through the COM Styles collection too
# (best effort - some hosts do not persist style renames).
$wb.Styles.Item("Normal").Name = "Standard"
